$d = $word.ActiveDocument

# 1. "min" -> "ደቂቃ" (single occurrence, time unit after "10")
$d.Content.Find.Execute("min", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ደቂቃ", 2)

# 2-4. "(Leave " / "BLANK" / " for the facilitators that will use it)" -> blanks (4 occurrences, replace all)
$d.Content.Find.Execute("(Leave ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "______________ ", 2)
$d.Content.Find.Execute("BLANK", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "_______", 2)
$d.Content.Find.Execute(" for the facilitators that will use it)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " ________________", 2)

# 5. "Ask the groups of students to " -> Amharic translation
$d.Content.Find.Execute("Ask the groups of students to ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "የተማሪዎቹን ቡዲኖች  ", 2)

# 6. "draw the square grids (you can draw and share a printed copy) and cut a few post-it papers the same size as the square grids." -> Amharic translation
$d.Content.Find.Execute("draw the square grids (you can draw and share a printed copy) and cut a few post-it papers the same size as the square grids.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ካሬ ፈርግርግ እነድሰሩ መጠየቅ (መሳልና ፕሪንት የሆነ ኮፒ ማጋራት ትችላላችሁ) እና የተወሰኑ በካሬው ፈርግርግ መጠን የወረቀት ልጥፎችን መቁረጥ", 2)
